# edit.ps1
# Implements commit "added report to all tests":
#  - Renames header D1 from TESTNAME to NAMEPARTICIPANT
#  - Adds 5 new header columns (O:S) for report-related fields
#  - Adds 3 new test rows (OpenViduReactTest, OpenViduVueTest, OpenViduHelloWordTest)
#  - Adds hyperlinks for the new rows' URL cells
#  - Widens the new columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D1 header: TESTNAME -> NAMEPARTICIPANT ---
$ws.Range("D1").Value = "NAMEPARTICIPANT"

# --- New header cells O1:S1 (copy header style from an existing header cell) ---
$ws.Range("N1").Copy($ws.Range("O1"))
$ws.Range("O1").Value = "idMainTitle"
$ws.Range("N1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = "xpathSessionName"
$ws.Range("N1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "XpathParticipantName"
$ws.Range("N1").Copy($ws.Range("R1"))
$ws.Range("R1").Value = "idHeaderStartPage"
$ws.Range("N1").Copy($ws.Range("S1"))
$ws.Range("S1").Value = "xpathLeaveButton"

# --- Row 4: OpenViduReactTest ---
$ws.Hyperlinks.Add($ws.Range("B4"), "http://localhost:3000/")
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = "OpenViduReactTest"
$ws.Range("B2").Copy($ws.Range("B4"))
$ws.Range("B4").Value = "http://localhost:3000/"
$ws.Range("A2").Copy($ws.Range("C4"))
$ws.Range("C4").Value = "TestSession"
$ws.Range("A2").Copy($ws.Range("D4"))
$ws.Range("D4").Value = "TestParticipant"
$ws.Range("A2").Copy($ws.Range("E4"))
$ws.Range("E4").Value = "//*[@id='join-dialog']/form/p[3]/input"
$ws.Range("A2").Copy($ws.Range("F4"))
$ws.Range("F4").Value = "buttonLeaveSession"
$ws.Range("A2").Copy($ws.Range("G4"))
$ws.Range("G4").Value = "/html/body/div/div/div/div[3]/div[2]/div/div/video"
$ws.Range("A2").Copy($ws.Range("H4"))
$ws.Range("H4").Value = "session-title"
$ws.Range("A2").Copy($ws.Range("I4"))
$ws.Range("I4").Value = "sessionId"
$ws.Range("A2").Copy($ws.Range("J4"))
$ws.Range("J4").Value = "local-video-undefined"
$ws.Range("A2").Copy($ws.Range("K4"))
$ws.Range("K4").Value = "userName"
$ws.Range("A2").Copy($ws.Range("L4"))
$ws.Range("L4").Value = "userName"
$ws.Range("A2").Copy($ws.Range("M4"))
$ws.Range("M4").Value = "//*[@id='main-video']/div/div/div/p"
$ws.Range("A2").Copy($ws.Range("O4"))
$ws.Range("O4").Value = "join"

# --- Row 5: OpenViduVueTest ---
$ws.Hyperlinks.Add($ws.Range("B5"), "http://localhost:8080/")
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "OpenViduVueTest"
$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("B5").Value = "http://localhost:8080/"
$ws.Range("A2").Copy($ws.Range("C5"))
$ws.Range("C5").Value = "TestSession"
$ws.Range("A2").Copy($ws.Range("D5"))
$ws.Range("D5").Value = "TestParticipant"
$ws.Range("A2").Copy($ws.Range("E5"))
$ws.Range("E5").Value = "//*[@id='join-dialog']/div/p[3]/button"
$ws.Range("A2").Copy($ws.Range("F5"))
$ws.Range("F5").Value = "buttonLeaveSession"
$ws.Range("A2").Copy($ws.Range("G5"))
$ws.Range("G5").Value = "/html/body/div/div/div[3]/div[2]/video"
$ws.Range("A2").Copy($ws.Range("H5"))
$ws.Range("H5").Value = "session-title"
$ws.Range("A2").Copy($ws.Range("J5"))
$ws.Range("J5").Value = "local-video-undefined"
$ws.Range("A2").Copy($ws.Range("M5"))
$ws.Range("M5").Value = "//*[@id='main-video']/div/div/p"
$ws.Range("A2").Copy($ws.Range("P5"))
$ws.Range("P5").Value = "//*[@id='join-dialog']/div/p[2]/input"
$ws.Range("A2").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = "//*[@id='join-dialog']/div/p[1]/input"
$ws.Range("A2").Copy($ws.Range("R5"))
$ws.Range("R5").Value = "img-div"

# --- Row 6: OpenViduHelloWordTest ---
$ws.Hyperlinks.Add($ws.Range("B6"), "http://localhost:8080/")
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "OpenViduHelloWordTest"
$ws.Range("B2").Copy($ws.Range("B6"))
$ws.Range("B6").Value = "http://localhost:8080/"
$ws.Range("A2").Copy($ws.Range("C6"))
$ws.Range("C6").Value = "TestSession"
$ws.Range("A2").Copy($ws.Range("E6"))
$ws.Range("E6").Value = "//*[@id='join']/form/p[2]/input"
$ws.Range("A2").Copy($ws.Range("G6"))
$ws.Range("G6").Value = "/html/body/div[2]/div/div[2]/video"
$ws.Range("A2").Copy($ws.Range("H6"))
$ws.Range("H6").Value = "session-header"
$ws.Range("A2").Copy($ws.Range("I6"))
$ws.Range("I6").Value = "sessionId"
$ws.Range("A2").Copy($ws.Range("J6"))
$ws.Range("J6").Value = "local-video-undefined"
$ws.Range("A2").Copy($ws.Range("S6"))
$ws.Range("S6").Value = "//*[@id='session']/input"

# --- Column widths for the new columns P:S (closest achievable to 17.13/19.5/15.88/15.63 chars) ---
$ws.Range("P1").ColumnWidth = 16.333333333333332
$ws.Range("Q1").ColumnWidth = 18.666666666666668
$ws.Range("R1").ColumnWidth = 15.0
$ws.Range("S1").ColumnWidth = 14.833333333333334

